$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.169743299484253
$ws.Range("B1").Value = 0.941416323184967
$ws.Range("C1").Value = 3.166744232177734
$ws.Range("D1").Value = 3.166005611419678
$ws.Range("E1").Value = 0.9354562759399414
